$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O (15th column) for the new
# "Server Session ID (Delegator)" variable. This shifts every existing
# column from O onward one place to the right (O->P, P->Q, ... W->X).
$ws.Columns.Item(15).Insert()

# Populate the new column's header / source / definition rows.
$ws.Range("O1").Value = "Server Session ID (Delegator)"
$ws.Range("O2").Value = "Text Scenario Output/TA1 server"
$ws.Range("O3").Value = "The session id used to get the delegator's alignment score from the TA1 server for the text scenarios"
# Row 4 ("Levels") intentionally left blank for this variable.

# Match the author's final selection/view state.
[void]$ws.Range("O4").Select()
